$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Helper: write a text value into a cell without Excel re-interpreting
# strings that look like percentages (e.g. "25.0%") as numbers. We
# stage the literal text in an always-unused helper cell (column J has
# no data anywhere in this sheet), force it to Text format there,
# then copy/paste-values (which does not re-parse the text) into the
# real destination, leaving the destination's original style/format
# untouched. The helper cell is fully cleared afterwards so it leaves
# no trace in the saved file.
# ------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $helper = $ws.Range("J1")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
}

# ------------------------------------------------------------------
# Row 2 (C1 / ANATOMY / session 1) - reorder "Recorded By" list
# ------------------------------------------------------------------
$ws.Range("G2").Value = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# ------------------------------------------------------------------
# Row 3 (C1 / ANATOMY / session 2)
# ------------------------------------------------------------------
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H3").Value = "106/221"

# ------------------------------------------------------------------
# Row 4 (C1 / ANATOMY / session 3)
# ------------------------------------------------------------------
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# ------------------------------------------------------------------
# Row 5 (C1 / ANATOMY / session 4) turns from "Recorded" into
# "Pending" - copy the formatting of an existing Pending row (row 6)
# onto row 5, then fix up the values.
# ------------------------------------------------------------------
$ws.Range("A6:I6").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = "0/221"
$ws.Range("I5").Value = "Pending"

# ------------------------------------------------------------------
# Statistics block next to the C1/ANATOMY group (column K/L)
# ------------------------------------------------------------------
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 41
Set-TextValue $ws.Range("L9") "23.2%"
Set-TextValue $ws.Range("L10") "41.6%"

# ------------------------------------------------------------------
# Row 13 (C1 / HISTOLOGY / session 1)
# ------------------------------------------------------------------
$ws.Range("G13").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# ------------------------------------------------------------------
# Group statistics rows 15 (C1) and 16 (C2)
# ------------------------------------------------------------------
$ws.Range("O15").Value = 7
$ws.Range("Q15").Value = 21
Set-TextValue $ws.Range("R15") "25.0%"
Set-TextValue $ws.Range("S15") "48.4%"

$ws.Range("O16").Value = 6
$ws.Range("P16").Value = 2
Set-TextValue $ws.Range("R16") "21.4%"
Set-TextValue $ws.Range("S16") "33.7%"

# ------------------------------------------------------------------
# Row 24 (C1 / PHYSIOLOGY / session 1)
# ------------------------------------------------------------------
$ws.Range("G24").Value = "marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("H24").Value = "89/221"

# ------------------------------------------------------------------
# Row 25 (C1 / PHYSIOLOGY / session 2)
# ------------------------------------------------------------------
$ws.Range("G25").Value = "abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("H25").Value = "82/221"

# ------------------------------------------------------------------
# Row 31 (C2 / ANATOMY / session 2)
# ------------------------------------------------------------------
$ws.Range("G31").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H31").Value = "42/246"

# ------------------------------------------------------------------
# Row 32 (C2 / ANATOMY / session 3)
# ------------------------------------------------------------------
$ws.Range("G32").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("H32").Value = "120/246"

# ------------------------------------------------------------------
# Row 33 (C2 / ANATOMY / session 4) turns from "Recorded" into
# "Not Recorded" - copy formatting of an existing Not Recorded row
# (row 30) onto row 33, then fix up the values.
# ------------------------------------------------------------------
$ws.Range("A30:I30").Copy()
$ws.Range("A33:I33").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G33").Value = ""
$ws.Range("H33").Value = "0/246"
$ws.Range("I33").Value = "Not Recorded"

# ------------------------------------------------------------------
# Row 41 (C2 / HISTOLOGY / session 1)
# ------------------------------------------------------------------
$ws.Range("G41").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("H41").Value = "70/246"

# ------------------------------------------------------------------
# Row 52 (C2 / PHYSIOLOGY / session 1)
# ------------------------------------------------------------------
$ws.Range("G52").Value = "marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("H52").Value = "91/246"

# ------------------------------------------------------------------
# Row 53 (C2 / PHYSIOLOGY / session 2)
# ------------------------------------------------------------------
$ws.Range("G53").Value = "abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("H53").Value = "28/246"
